$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells are plain text (t="inlineStr"), including values that look
# numeric ("312.08", "35.30", ...). Assigning a bare numeric-looking string
# via .Value lets Excel coerce it to a real number (losing trailing zeros /
# introducing float noise), so each value is prefixed with a literal leading
# apostrophe - Excel's "force text" entry convention; the apostrophe itself
# is not stored - and .Style is reset to "Normal" afterwards so no stray
# quote-prefix / number-format style stays attached to the cell.

$ws.Range("D2").Value = "`'44.528.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "`'  +3.95%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "`'2.430.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "`'  +2.98%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "`'  -0.17%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "`'312.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "`'  +3.63%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "`'101.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "`'  +7.08%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "`'  +2.24%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "`'  -0.07%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "`'0.511"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "`'  +5.65%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "`'35.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "`'  +4.53%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "`'0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "`'  +2.03%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "`'0.124"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "`'  +1.56%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "`'18.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "`'  +3.42%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "`'  +3.70%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "`'2.810.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "`'  +2.79%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "`'2.397.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "`'  +1.65%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "`'0.839"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "`'  +5.70%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "`'44.434.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "`'  +3.82%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "`'12.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "`'  +3.58%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "`'6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "`'  +2.36%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "`'  +2.78%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "`'68.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "`'  +1.65%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "`'241.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "`'  +2.73%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "`'2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "`'  +4.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "`'  +2.79%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "`'  -0.02%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "`'25.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "`'  +2.55%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E29").Value = "`'  +4.86%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "`'33.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "`'  +6.50%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "`'48.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "`'  +1.49%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "`'0.122"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "`'  +17.55%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "`'19.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "`'  +13.03%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "`'  +3.81%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "`'  +0.15%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "`'0.0764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "`'  +5.61%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "`'  +3.17%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "`'4.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "`'  +4.32%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "`'2.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "`'  +4.88%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "`'126.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "`'  +4.85%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "`'Stellar"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "`'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "`'0.109"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "`'  +1.40%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "`'WEMIXToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "`'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "`'2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "`'  -5.44%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "`'21.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "`'  +0.68%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "`'  +4.03%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "`'1.946.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "`'  +0.87%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "`'  +2.28%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "`'2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "`'  +8.93%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "`'9.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "`'  +7.09%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "`'  +12.30%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "`'53.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "`'  +4.04%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "`'73.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "`'  +2.66%  "
$ws.Range("E51").Style = "Normal"
